$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()
$ws.Range("A545:XFD545").Select()
"done"
